# edit.ps1 - Applies the "NOTAS.docx" diff: appends a large block of new
# paragraphs (dated note + TODO checklist + "TEMAS A REVISAR" section) to the
# end of the document body, right before the trailing <w:sectPr>.
#
# Strategy: build the exact OOXML for the new paragraphs and insert it with
# Range.InsertXML at a *freshly created* Range anchored on the document's
# current end-of-content position. (Reusing a Paragraph's own .Range object
# after .Collapse() was found to insert at the wrong offset in this runtime;
# creating a brand-new $d.Range(pos, pos) each time is reliable.)

$d = $word.ActiveDocument

$newParagraphsXml = '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>14/9/23</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Actualmente se genero la variable global del personaje. Aun no acomodo la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pagina</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> por lo que no se verifico el funcionamiento correcto de la variable. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Revisar que funcione correctamente la variable </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PJ_active</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Que muestre nombre de personaje.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Muestre nivel de personaje</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Muestre vida de personaje</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Mostrar imagen de personaje</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Hacer que </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>system_m</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> sea = 0, agregarlo a city.html y acomodar estilo de escenario.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Revisar sistema de inventario y los pasos de su funcionamiento.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Modificar base de datos y código para guardar los nuevos datos. (dinero e inventario).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">acomodar </w:t></w:r><w:r><w:t xml:space="preserve">menú de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>city</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para </w:t></w:r><w:r><w:t>que funcione correctamente.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Agregar ventanas de inventario, equipo y estadísticas.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Que el inventario muestre los objetos y se pueda organizar. Que pueda inspeccionar objeto para mostrar detalles.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Equipo: que muestre la ventana de objetos equipados, que muestre detalles de objetos y pueda cambiarlos.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Estadísticas: que muestre estadísticas del personaje y pueda subirlas (modificar personaje para agregar puntos de estadísticas.)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Verifica si tiene puntos de estadísticas para mostrar icono + (&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>img</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&gt;color </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>gris ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> o &lt;a&gt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>img</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;color blanco si hay puntos)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Ejecutar función para subir estadística correspondiente y actualizar estadísticas.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>PARA ESTE PUNTO YA MODIFICASTE SEGURAMENTE LAS COSAS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Agregar a las opciones la función de guardar.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Revisar función cargar partida para ver si toma los cambios y como guardaron las cosas modificadas en el juego.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>(modificaciones de ser necesario)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Acomodar estilo de juego (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>FRAMEWORK?</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Cambio de escenario =&gt; tienda</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>TEMAS A REVISAR</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Cuando se carga city.html hay momentos donde toma los datos de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PJ_active</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y otros donde dice que es </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>null</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, revisar. (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">ver si se puede colocar que cargue todos los datos antes o que cargue la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pagina</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y después modifique los datos?</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Revisar el menú. Separar el sistema de menú en archivos apartes. Tanto </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>html</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> como </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>javascript</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> y CSS. Revisar que cuando cargue el sistema de menú muestre el inicio. Hasta ahora no lo muestra.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Revisar cálculos de estadísticas personaje.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$openXml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' + $newParagraphsXml + '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$endPos = $d.Content.End
$insertionRange = $d.Range($endPos, $endPos)
$insertionRange.InsertXML($openXml)

Write-Output "Inserted new paragraphs; document now has $($d.Paragraphs.Count) paragraphs."
